$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 8, shifting rows 8:23 down to 9:24.
$ws.Rows.Item(8).Insert()

# Fill the newly inserted row 8 with data.
$ws.Range("A8").Value = -0.08575
$ws.Range("B8").Value = 0.08575
$ws.Range("C8").Value = 100
$ws.Range("D8").Value = 100
$ws.Range("E8").Value = 0.0148212
$ws.Range("F8").Value = 2
$ws.Range("G8").Value = 12

# New header columns I1:L1
$ws.Range("I1").Value = "modelRTR"
$ws.Range("J1").Value = "modelLR"
$ws.Range("K1").Value = "modelXGB"
$ws.Range("L1").Value = "score"

# New data rows 2-4 for columns I:L
$ws.Range("I2").Value = 0.6
$ws.Range("J2").Value = 0.3
$ws.Range("K2").Value = 0.1
$ws.Range("L2").Value = 0.0113736

$ws.Range("I3").Value = 0.7
$ws.Range("J3").Value = 0.3
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0.0139797

$ws.Range("I4").Value = 0.5
$ws.Range("J4").Value = 0.2
$ws.Range("K4").Value = 0.3

# Update the active selection to match target state
$ws.Range("P6").Select()
